$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sheet name and header title to reflect new "through" date ---
$ws.Name = "Through 2022-03-22"
$ws.Range("B1").Value = "March 2022 (through March 22)"

# --- Simple single-cell updates (new data for 2022-03-30 added to existing rows) ---
$ws.Range("T4").Value = 4

$ws.Range("B7").Value = 2
$ws.Range("E7").Value = 4

$ws.Range("E10").Value = 1

$ws.Range("D11").Value = 13

$ws.Range("B14").Value = 2

$ws.Range("N15").Value = 2

$ws.Range("N18").Value = 2

$ws.Range("N48").Value = 1

$ws.Range("B75").Value = 1

$ws.Range("B77").Value = 3

# --- Rows 22-26 (Chatham, Bridgeport, Wicker Park, Washington Park, Grand Crossing) ---
# Adding the new day's counts shifted each neighborhood's rank by one row, so the
# entire block of rows is rewritten with its new label + values.

# Row 22 -> Bridgeport
$ws.Range("A22").Value = "Bridgeport"
$ws.Range("B22:Y22").ClearContents()
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 4
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("I22").Value = 1
$ws.Range("N22").Value = 1
$ws.Range("S22").Value = 1

# Row 23 -> Wicker Park
$ws.Range("A23").Value = "Wicker Park"
$ws.Range("B23:Y23").ClearContents()
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 2
$ws.Range("L23").Value = 1
$ws.Range("R23").Value = 1

# Row 24 -> Washington Park
$ws.Range("A24").Value = "Washington Park"
$ws.Range("B24:Y24").ClearContents()
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = 2
$ws.Range("E24").Value = 1
$ws.Range("J24").Value = 2
$ws.Range("L24").Value = 1
$ws.Range("O24").Value = 1
$ws.Range("P24").Value = 4
$ws.Range("R24").Value = 3
$ws.Range("U24").Value = 1
$ws.Range("V24").Value = 1

# Row 25 -> Grand Crossing
$ws.Range("A25").Value = "Grand Crossing"
$ws.Range("B25:Y25").ClearContents()
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 2
$ws.Range("J25").Value = 1
$ws.Range("L25").Value = 3
$ws.Range("M25").Value = 1
$ws.Range("N25").Value = 2
$ws.Range("O25").Value = 2
$ws.Range("P25").Value = 1
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = 2
$ws.Range("S25").Value = 3
$ws.Range("T25").Value = 2
$ws.Range("U25").Value = 1
$ws.Range("V25").Value = 2
$ws.Range("Y25").Value = 1

# Row 26 -> Chatham
$ws.Range("A26").Value = "Chatham"
$ws.Range("B26:Y26").ClearContents()
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 4
$ws.Range("L26").Value = 1
$ws.Range("N26").Value = 1
$ws.Range("P26").Value = 1
$ws.Range("U26").Value = 2
$ws.Range("W26").Value = 1
$ws.Range("Y26").Value = 1
